# Auto-generated script applying cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.658.19"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.06%  "
$ws.Cells.Item(3, 4).Value = "'1.597.50"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.20%  "
$ws.Cells.Item(5, 4).Value = "'211.48"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.30%  "
$ws.Cells.Item(6, 4).Value = "'0.513"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.04%  "
$ws.Cells.Item(7, 5).Value = "  +0.14%  "
$ws.Cells.Item(8, 4).Value = "'0.0618"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.12%  "
$ws.Cells.Item(9, 5).Value = "  +0.65%  "
$ws.Cells.Item(10, 4).Value = "'19.49"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -0.79%  "
$ws.Cells.Item(11, 4).Value = "'0.0841"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.32%  "
$ws.Cells.Item(12, 4).Value = "'1.822.39"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.02%  "
$ws.Cells.Item(13, 4).Value = "'1.597.03"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.02%  "
$ws.Cells.Item(14, 5).Value = "  +0.21%  "
$ws.Cells.Item(15, 5).Value = "  +0.38%  "
$ws.Cells.Item(16, 4).Value = "'65.01"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.31%  "
$ws.Cells.Item(17, 4).Value = "'26.640.21"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.04%  "
$ws.Cells.Item(18, 5).Value = "  +1.29%  "
$ws.Cells.Item(19, 5).Value = "  +0.16%  "
$ws.Cells.Item(20, 4).Value = "'208.78"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.12%  "
$ws.Cells.Item(21, 4).Value = "'7.03"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +3.93%  "
$ws.Cells.Item(22, 4).Value = "'4.28"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.49%  "
$ws.Cells.Item(23, 4).Value = "'2.35"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +2.50%  "
$ws.Cells.Item(24, 5).Value = "  +1.18%  "
$ws.Cells.Item(25, 4).Value = "'144.12"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -1.21%  "
$ws.Cells.Item(26, 5).Value = "  +0.21%  "
$ws.Cells.Item(27, 4).Value = "'7.12"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -1.77%  "
$ws.Cells.Item(28, 5).Value = "  -0.92%  "
$ws.Cells.Item(29, 5).Value = "  +0.10%  "
$ws.Cells.Item(30, 5).Value = "  +1.62%  "
$ws.Cells.Item(31, 5).Value = "  +0.24%  "
$ws.Cells.Item(32, 4).Value = "'3.23"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.24%  "
$ws.Cells.Item(33, 4).Value = "'2.94"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.82%  "
$ws.Cells.Item(34, 4).Value = "'1.290.20"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.38%  "
$ws.Cells.Item(35, 4).Value = "'0.619"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -7.05%  "
$ws.Cells.Item(36, 5).Value = "  +0.54%  "
$ws.Cells.Item(37, 4).Value = "'1.48"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.11%  "
$ws.Cells.Item(38, 5).Value = "  -0.18%  "
$ws.Cells.Item(39, 4).Value = "'0.831"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.25%  "
$ws.Cells.Item(40, 4).Value = "'1.03"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +15.30%  "
$ws.Cells.Item(41, 5).Value = "  +1.44%  "
$ws.Cells.Item(42, 5).Value = "  -0.45%  "
$ws.Cells.Item(43, 5).Value = "  -0.30%  "
$ws.Cells.Item(44, 4).Value = "'63.30"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.83%  "
$ws.Cells.Item(45, 4).Value = "'1.733.98"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.08%  "
$ws.Cells.Item(46, 4).Value = "'91.13"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.24%  "
$ws.Cells.Item(47, 5).Value = "  -2.68%  "
$ws.Cells.Item(48, 2).Value = "Algorand"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(48, 4).Value = "'0.101"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.26%  "
$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(49, 4).Value = "'0.0509"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.97%  "
$ws.Cells.Item(50, 2).Value = "USDD"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Cells.Item(50, 4).Value = "'1.00"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.04%  "
$ws.Cells.Item(51, 2).Value = "EnergySwap"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(51, 4).Value = "'7.36"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.82%  "
